$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries (rows 26-29): preposition words
$ws.Range("B26").Value = "ał "
$ws.Range("B26").WrapText = $true
$ws.Range("C26").Value = "at"
$ws.Range("E26").Value = "preposition"

$ws.Range("B27").Value = "łaxi"
$ws.Range("C27").Value = "under"
$ws.Range("E27").Value = "preposition"

$ws.Range("B28").Value = "lax̱'oi"
$ws.Range("C28").Value = "on"
$ws.Range("E28").Value = "preposition"

$ws.Range("B29").Value = "ts'm"
$ws.Range("C29").Value = "in"
$ws.Range("E29").Value = "preposition"

# Update selection to match the author's final cursor position
$ws.Range("B21").Select()
